$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the VNET Address Space column (B2:B6) entirely
$ws.Range("B2:B6").ClearContents()

# Row 3: remove Subnet Address Space, NSG, UDR (D3, F3, G3)
$ws.Range("D3").ClearContents()
$ws.Range("F3:G3").Clear()

# Row 4: add new "Gateway" entry in Subnet Name column (C4)
$ws.Range("C4").Value = "Gateway"

# Update selection to match the new active cell
$ws.Range("C4").Select()
